# Updated symbol list on Mon Jan  9 06:43:23 UTC 2023 with GitHub Actions
#
# Refresh the crypto price / 1h-volume columns (D, E) for this pull of the
# symbol list. All of these columns are stored as plain text in the sheet
# (no numeric formatting), so each new value is written with a leading
# apostrophe to force Excel to keep it as literal text instead of
# re-interpreting it as a number/percentage (which would silently change
# the stored representation, e.g. dropping a trailing zero like "0.1450").

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value  = "'278.35"

$ws.Range("D3").Value  = "'27.26"
$ws.Range("E3").Value  = "'1.65%"

$ws.Range("D4").Value  = "'4.847"
$ws.Range("E4").Value  = "'3.35%"

$ws.Range("D5").Value  = "'0.06248"
$ws.Range("E5").Value  = "'0.51%"

$ws.Range("D6").Value  = "'6.897"
$ws.Range("E6").Value  = "'2.12%"

$ws.Range("D7").Value  = "'0.8784"
$ws.Range("E7").Value  = "'3.17%"

$ws.Range("D8").Value  = "'0.9435"
$ws.Range("E8").Value  = "'3.21%"

$ws.Range("D9").Value  = "'0.1450"
$ws.Range("E9").Value  = "'3.47%"

$ws.Range("E10").Value = "'6.47%"

$ws.Range("D11").Value = "'0.07319"
$ws.Range("E11").Value = "'3.21%"

$ws.Range("D12").Value = "'0.03159"
$ws.Range("E12").Value = "'1.87%"

$ws.Range("E13").Value = "'0.06%"

$ws.Range("D14").Value = "'0.001562"
$ws.Range("E14").Value = "'1.93%"

$ws.Range("D15").Value = "'0.0006272"
$ws.Range("E15").Value = "'1.71%"

$ws.Range("D16").Value = "'0.006114"
$ws.Range("E16").Value = "'0.24%"

$ws.Range("D17").Value = "'3.452"
$ws.Range("E17").Value = "'0.28%"

$ws.Range("E18").Value = "'2.70%"

$ws.Range("E19").Value = "'5.59%"

$ws.Range("D21").Value = "'0.1311"
$ws.Range("E21").Value = "'0.03%"

$ws.Range("D22").Value = "'3.845"
$ws.Range("E22").Value = "'-5.82%"

$ws.Range("E23").Value = "'1.98%"

$ws.Range("E24").Value = "'-2.28%"

$ws.Range("D25").Value = "'0.004276"
$ws.Range("E25").Value = "'4.84%"

$ws.Range("D27").Value = "'0.0001690"
$ws.Range("E27").Value = "'3.05%"

$ws.Range("D40").Value = "'0.04028"
$ws.Range("E40").Value = "'1.73%"

$ws.Range("D41").Value = "'0.006483"
$ws.Range("E41").Value = "'56.89%"

$ws.Range("D42").Value = "'0.1153"
$ws.Range("E42").Value = "'3.49%"

$ws.Range("E43").Value = "'-4.67%"

$ws.Range("D44").Value = "'0.01199"
$ws.Range("E44").Value = "'-13.69%"

$ws.Range("D45").Value = "'0.00005093"
$ws.Range("E45").Value = "'-1.35%"

$ws.Range("D47").Value = "'2.370"
$ws.Range("E47").Value = "'854.66%"
